# ============================================================================
# Reproduces the "Add files via upload" commit:
#   - Demands sheet becomes the active tab (Ratings loses tabSelected)
#   - Demands sheet gains a split "Practitioner / Researcher" breakdown
#     (two new narrow columns C,D; old Responses column C becomes Total
#     in column E) plus a second header row.
#   - A blank spacer row is inserted under the main header.
#   - Styles: new borders (split top/bottom box around the two-row header,
#     left-only border framing "Responses"/"Total"), a tweaked italic-grey
#     font, and reuse of existing header/zebra fills.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ratings = $wb.Worksheets.Item("Ratings")
$ws = $wb.Worksheets.Item("Demands")

# ----------------------------------------------------------------------
# 1. Structural edits: insert the 2 "Practitioner/Researcher" columns
#    before the old column C, and insert one new row below the header
#    row for the new sub-header (Practitioner/Researcher/Total).
# ----------------------------------------------------------------------
$ws.Range("C1:D1").EntireColumn.Insert()
$ws.Range("A2").EntireRow.Insert()

# ----------------------------------------------------------------------
# 2. Header text (row 1 + new row 2)
# ----------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Demands"
$ws.Cells.Item(1,2).Value = "Statements"
$ws.Cells.Item(1,3).Value = "Responses"

$ws.Cells.Item(2,3).Value = "Practitioner"
$ws.Cells.Item(2,4).Value = "Researcher"
$ws.Cells.Item(2,5).Value = "Total"

# ----------------------------------------------------------------------
# 3. Fill in Practitioner / Researcher counts for each data row
#    (old "Responses" values already sit in column E from the column
#    insert - they become the "Total" column, so only C/D are new.)
# ----------------------------------------------------------------------
$data = @(
  @{ Row=3;  C=23; D=20 },
  @{ Row=5;  C=16; D=21 },
  @{ Row=6;  C=9;  D=7  },
  @{ Row=7;  C=4;  D=1  },
  @{ Row=9;  C=8;  D=14 },
  @{ Row=11; C=3;  D=16 },
  @{ Row=12; C=2;  D=12 },
  @{ Row=13; C=1;  D=4  },
  @{ Row=14; C=2;  D=8  },
  @{ Row=15; C=0;  D=2  },
  @{ Row=17; C=10; D=2  },
  @{ Row=19; C=3;  D=13 },
  @{ Row=20; C=0;  D=9  }
)
foreach ($d in $data) {
  $ws.Cells.Item($d.Row, 3).Value = $d.C
  $ws.Cells.Item($d.Row, 4).Value = $d.D
}

# ----------------------------------------------------------------------
# 4. Merged cells matching the new layout
# ----------------------------------------------------------------------
$ws.Range("A1:A2").Merge()
$ws.Range("B1:B2").Merge()
$ws.Range("C1:E1").Merge()

# ----------------------------------------------------------------------
# 5. Column widths (2 new narrow cols C:D, 1 new wide col E)
# ----------------------------------------------------------------------
$ws.Range("C1").EntireColumn.ColumnWidth = 10.83203125
$ws.Range("D1").EntireColumn.ColumnWidth = 10.83203125
$ws.Range("E1").EntireColumn.ColumnWidth = 20.83203125

# Row heights for the two header rows
$ws.Rows.Item(1).RowHeight = 28
$ws.Rows.Item(2).RowHeight = 28

# ----------------------------------------------------------------------
# 6. Formatting - reuse the existing "big header" look (bold 14pt, green
#    fill, thin border) for A1/B1 but give the 2-row block a split
#    border (top-half has no bottom line, bottom-half has no top line
#    so the merged block reads as a single boxed cell).
# ----------------------------------------------------------------------
foreach ($addr in @("A1","A2","B1","B2")) {
  $c = $ws.Range($addr)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 14
  $c.Font.Bold = $true
  $c.Font.Italic = $false
  $c.Font.ThemeColor = 1
  $c.Interior.ThemeColor = 10
  $c.Interior.TintAndShade = 0.79998168889431442
  $c.HorizontalAlignment = -4108
  $c.VerticalAlignment = -4108
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(10).LineStyle = 1
}
foreach ($addr in @("A1","B1")) {
  $c = $ws.Range($addr)
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(9).LineStyle = 0
}
foreach ($addr in @("A2","B2")) {
  $c = $ws.Range($addr)
  $c.Borders.Item(8).LineStyle = 0
  $c.Borders.Item(9).LineStyle = 1
}

# "Responses" header (C1:E1 merged) - same bold/green look, but only a
# left border is drawn (the sub-header row underneath closes the box).
foreach ($addr in @("C1","D1","E1")) {
  $c = $ws.Range($addr)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 14
  $c.Font.Bold = $true
  $c.Font.Italic = $false
  $c.Font.ThemeColor = 1
  $c.Interior.ThemeColor = 10
  $c.Interior.TintAndShade = 0.79998168889431442
  $c.HorizontalAlignment = -4108
  $c.VerticalAlignment = -4108
  $c.Borders.Item(7).LineStyle = 0
  $c.Borders.Item(8).LineStyle = 0
  $c.Borders.Item(9).LineStyle = 0
  $c.Borders.Item(10).LineStyle = 0
}
$ws.Range("C1").Borders.Item(7).LineStyle = 1

# Practitioner / Researcher / Total sub-header (row 2, C:E) - normal
# weight green header cell, fully boxed.
foreach ($addr in @("C2","D2","E2")) {
  $c = $ws.Range($addr)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 12
  $c.Font.Bold = $false
  $c.Font.Italic = $false
  $c.Font.ThemeColor = 1
  $c.Interior.ThemeColor = 10
  $c.Interior.TintAndShade = 0
  $c.HorizontalAlignment = -4108
  $c.VerticalAlignment = -4108
  $c.Borders.Item(7).LineStyle = 1
  $c.Borders.Item(8).LineStyle = 1
  $c.Borders.Item(9).LineStyle = 1
  $c.Borders.Item(10).LineStyle = 1
}

Write-Output "structure-and-headers-done"
